$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status": Insufficient Data row totals bump 301 -> 302 ---
$trends = $wb.Worksheets.Item("Trends Status")
$trends.Range("B8").Value = 302
$trends.Range("C8").Value = 302

# --- Sheet "Priority Status": counts updated ---
$priority = $wb.Worksheets.Item("Priority Status")
$priority.Range("B2").Value = 103
$priority.Range("B3").Value = 286
$priority.Range("B4").Value = 554

# --- Sheet "Species qualification": label + count updated ---
$qual = $wb.Worksheets.Item("Species qualification")
$qual.Range("A2").Value = "SoIB Assessment"
$qual.Range("B2").Value = 302

# --- Sheet "High Priority break-up": duplicate as a new sheet, then ---
# --- repurpose the original as the "Interannual update" breakdown    ---
$src = $wb.Worksheets.Item("High Priority break-up")
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Major update - High Priority "
$src.Name = "Interannual update - High Pri"

# Write the (moved-down) IUCN row first so we don't disturb row 2 values
# while they're being read/written.
$src.Range("A3").Value = "IUCN"
$src.Range("B3").Value = 11
$src.Range("C3").Value = 10.7
$src.Range("D3").Value = 7
$src.Range("E3").Value = 7.1

# Now set the new "Trend New" row in row 2.
$src.Range("A2").Value = "Trend New"
$src.Range("B2").Value = 92
$src.Range("C2").Value = 89.3
$src.Range("D2").Value = 92
$src.Range("E2").Value = 92.90000000000001

# Restore original active sheet/tab.
$wb.Worksheets.Item("Trends Status").Select()
